# Updated cryptos list on Fri Aug 23 23:52:59 UTC 2024 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest scraped snapshot. Values are stored as plain text (inline strings),
# matching the existing sheet convention, so numeric-looking prices are
# assigned with a leading apostrophe to keep Excel from reinterpreting them
# as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.030.70"
$ws.Range("E2").Value = "  +6.20%  "
$ws.Range("D3").Value = "2.732.81"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'591.47"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'153.16"
$ws.Range("E6").Value = "  +6.88%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "2.765.34"
$ws.Range("E9").Value = "  +5.44%  "
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("E11").Value = "  +7.18%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "3.221.59"
$ws.Range("E14").Value = "  +4.80%  "
$ws.Range("D15").Value = "'26.59"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").Value = "63.887.30"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").Value = "'0.0000153"
$ws.Range("E17").Value = "  +8.76%  "
$ws.Range("D18").Value = "2.760.85"
$ws.Range("D19").Value = "'12.06"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("D21").Value = "'365.47"
$ws.Range("E21").Value = "  +5.42%  "
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'0.539"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'66.19"
$ws.Range("E25").Value = "  +3.94%  "
$ws.Range("E26").Value = "  +5.22%  "
$ws.Range("E27").Value = "  +7.92%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "0.0₃0911"
$ws.Range("E29").Value = "  +13.84%  "
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("E31").Value = "  +9.12%  "
$ws.Range("D32").Value = "'172.96"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +17.92%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("D36").Value = "'4.84"
$ws.Range("E36").Value = "  +12.79%  "
$ws.Range("D37").Value = "'1.44"
$ws.Range("E37").Value = "  +10.26%  "
$ws.Range("E38").Value = "  +8.86%  "
$ws.Range("D39").Value = "'1.01"
$ws.Range("E39").Value = "  +18.96%  "
$ws.Range("D40").Value = "'348.29"
$ws.Range("E40").Value = "  +8.88%  "
$ws.Range("E41").Value = "  +7.06%  "
$ws.Range("D42").Value = "'39.01"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("E43").Value = "  +10.43%  "
$ws.Range("D44").Value = "'5.59"
$ws.Range("E44").Value = "  +10.71%  "
$ws.Range("D45").Value = "'143.48"
$ws.Range("E45").Value = "  +5.69%  "
$ws.Range("E46").Value = "  +10.81%  "
$ws.Range("D47").Value = "'0.0593"
$ws.Range("E47").Value = "  +7.18%  "
$ws.Range("E48").Value = "  +6.47%  "
$ws.Range("E49").Value = "  +7.00%  "
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").Value = "2.176.67"
$ws.Range("E51").Value = "  +7.45%  "
